# Auto-generated edit script applying cryptocurrency price/volume updates
# Commit: Updated cryptos list on Mon May 22 11:07:50 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.016.53"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.829.61"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "311.78"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "0.4653"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").Value = "0.07384"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "0.8784"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07919"
$ws.Range("E11").Value = "  +8.09%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "20.01"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").Value = "1.760.84"
$ws.Range("E13").Value = "  -8.71%  "
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "6.592"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "92.10"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.000008887"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "27.051.37"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "5.169"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "2.061.89"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "152.66"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "1.835"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "18.29"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "2.097"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "5.136"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "115.58"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "0.08880"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "2.985"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "0.7295"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("D36").Value = "2.489"
$ws.Range("E36").Value = "  +3.35%  "
$ws.Range("D37").Value = "1.078"
$ws.Range("D38").Value = "0.01957"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "0.05252"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "7.343"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").Value = "0.5219"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "0.8588"
$ws.Range("E44").Value = "  -14.99%  "
$ws.Range("D45").Value = "8.241"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("D46").Value = "0.4862"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "10.24"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").Value = "102.88"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "1.627"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "0.06249"
$ws.Range("E51").Value = "  -0.87%  "
